$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 35.42516366666666
$ws.Range("H2").Value = 106.275491
$ws.Range("I2").Value = 0.00832770193000585
$ws.Range("J2").Value = 0.008327701930005852
$ws.Range("M2").Value = 0.5001966666666666
$ws.Range("N2").Value = 1.50059
$ws.Range("O2").Value = 0.03894027965151046
$ws.Range("P2").Value = 0.03894027965151046
$ws.Range("Q2").Value = 17.71954878218778
$ws.Range("R2").Value = 159.47593903969
$ws.Range("S2").Value = 0.0003242830420088512
$ws.Range("T2").Value = 0.0003242830420088513
$ws.Range("G3").Value = 35.42516366666666
$ws.Range("H3").Value = 106.275491
$ws.Range("I3").Value = 0.00832770193000585
$ws.Range("J3").Value = 0.008327701930005852
$ws.Range("O3").Value = 0.7732779360092192
$ws.Range("P3").Value = 0.7732779360092191
$ws.Range("Q3").Value = 351.8756473227675
$ws.Range("R3").Value = 3166.880825904907
$ws.Range("S3").Value = 0.006439628160134915
$ws.Range("T3").Value = 0.006439628160134916
$ws.Range("G4").Value = 35.42516366666666
$ws.Range("H4").Value = 106.275491
$ws.Range("I4").Value = 0.00832770193000585
$ws.Range("J4").Value = 0.008327701930005852
$ws.Range("M4").Value = 2.334238666666666
$ws.Range("N4").Value = 7.002715999999999
$ws.Range("O4").Value = 0.1817203362411497
$ws.Range("P4").Value = 0.1817203362411496
$ws.Range("Q4").Value = 82.69078680372843
$ws.Range("R4").Value = 744.2170812335559
$ws.Range("S4").Value = 0.001513312794836734
$ws.Range("T4").Value = 0.001513312794836734
$ws.Range("G5").Value = 35.42516366666666
$ws.Range("H5").Value = 106.275491
$ws.Range("I5").Value = 0.00832770193000585
$ws.Range("J5").Value = 0.008327701930005852
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.07786066666666666
$ws.Range("N5").Value = 0.233582
$ws.Range("O5").Value = 0.006061448098120818
$ws.Range("P5").Value = 0.006061448098120817
$ws.Range("Q5").Value = 2.758226859862444
$ws.Range("R5").Value = 24.82404173876199
$ws.Range("S5").Value = 0.00005047793302535102
$ws.Range("T5").Value = 0.00005047793302535102
$ws.Range("I6").Value = 0.01070182047907406
$ws.Range("J6").Value = 0.01070182047907406
$ws.Range("M6").Value = 0.5001966666666666
$ws.Range("N6").Value = 1.50059
$ws.Range("O6").Value = 0.03894027965151046
$ws.Range("P6").Value = 0.03894027965151046
$ws.Range("Q6").Value = 22.77115963455667
$ws.Range("R6").Value = 204.94043671101
$ws.Range("S6").Value = 0.0004167318822354054
$ws.Range("T6").Value = 0.0004167318822354054
$ws.Range("I7").Value = 0.01070182047907406
$ws.Range("J7").Value = 0.01070182047907406
$ws.Range("O7").Value = 0.7732779360092192
$ws.Range("P7").Value = 0.7732779360092191
$ws.Range("S7").Value = 0.008275481651599579
$ws.Range("T7").Value = 0.008275481651599578
$ws.Range("I8").Value = 0.01070182047907406
$ws.Range("J8").Value = 0.01070182047907406
$ws.Range("M8").Value = 2.334238666666666
$ws.Range("N8").Value = 7.002715999999999
$ws.Range("O8").Value = 0.1817203362411497
$ws.Range("P8").Value = 0.1817203362411496
$ws.Range("Q8").Value = 106.2648451019027
$ws.Range("R8").Value = 956.383605917124
$ws.Range("S8").Value = 0.001944738415849759
$ws.Range("T8").Value = 0.001944738415849758
$ws.Range("I9").Value = 0.01070182047907406
$ws.Range("J9").Value = 0.01070182047907406
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.07786066666666666
$ws.Range("N9").Value = 0.233582
$ws.Range("O9").Value = 0.006061448098120818
$ws.Range("P9").Value = 0.006061448098120817
$ws.Range("Q9").Value = 3.544561145788667
$ws.Range("R9").Value = 31.901050312098
$ws.Range("S9").Value = 0.00006486852938931386
$ws.Range("T9").Value = 0.00006486852938931384
$ws.Range("G10").Value = 51.06824600000001
$ws.Range("H10").Value = 153.204738
$ws.Range("I10").Value = 0.01200505761322374
$ws.Range("J10").Value = 0.01200505761322374
$ws.Range("M10").Value = 0.5001966666666666
$ws.Range("N10").Value = 1.50059
$ws.Range("O10").Value = 0.03894027965151046
$ws.Range("P10").Value = 0.03894027965151046
$ws.Range("Q10").Value = 25.54416642171334
$ws.Range("R10").Value = 229.89749779542
$ws.Range("S10").Value = 0.0004674803006914272
$ws.Range("T10").Value = 0.0004674803006914272
$ws.Range("G11").Value = 51.06824600000001
$ws.Range("H11").Value = 153.204738
$ws.Range("I11").Value = 0.01200505761322374
$ws.Range("J11").Value = 0.01200505761322374
$ws.Range("O11").Value = 0.7732779360092192
$ws.Range("P11").Value = 0.7732779360092191
$ws.Range("Q11").Value = 507.2572786952827
$ws.Range("R11").Value = 4565.315508257544
$ws.Range("S11").Value = 0.009283246172825416
$ws.Range("T11").Value = 0.009283246172825416
$ws.Range("G12").Value = 51.06824600000001
$ws.Range("H12").Value = 153.204738
$ws.Range("I12").Value = 0.01200505761322374
$ws.Range("J12").Value = 0.01200505761322374
$ws.Range("M12").Value = 2.334238666666666
$ws.Range("N12").Value = 7.002715999999999
$ws.Range("O12").Value = 0.1817203362411497
$ws.Range("P12").Value = 0.1817203362411496
$ws.Range("Q12").Value = 119.2054744520453
$ws.Range("R12").Value = 1072.849270068408
$ws.Range("S12").Value = 0.002181563106069392
$ws.Range("T12").Value = 0.002181563106069392
$ws.Range("G13").Value = 51.06824600000001
$ws.Range("H13").Value = 153.204738
$ws.Range("I13").Value = 0.01200505761322374
$ws.Range("J13").Value = 0.01200505761322374
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.07786066666666666
$ws.Range("N13").Value = 0.233582
$ws.Range("O13").Value = 0.006061448098120818
$ws.Range("P13").Value = 0.006061448098120817
$ws.Range("Q13").Value = 3.976207679057334
$ws.Range("R13").Value = 35.785869111516
$ws.Range("S13").Value = 0.00007276803363750587
$ws.Range("T13").Value = 0.00007276803363750587
$ws.Range("G14").Value = 4121.876464666667
$ws.Range("H14").Value = 12365.629394
$ws.Range("I14").Value = 0.9689654199776964
$ws.Range("J14").Value = 0.9689654199776964
$ws.Range("M14").Value = 0.5001966666666666
$ws.Range("N14").Value = 1.50059
$ws.Range("O14").Value = 0.03894027965151046
$ws.Range("P14").Value = 0.03894027965151046
$ws.Range("Q14").Value = 2061.748868038051
$ws.Range("R14").Value = 18555.73981234246
$ws.Range("S14").Value = 0.03773178442657478
$ws.Range("T14").Value = 0.03773178442657478
$ws.Range("G15").Value = 4121.876464666667
$ws.Range("H15").Value = 12365.629394
$ws.Range("I15").Value = 0.9689654199776964
$ws.Range("J15").Value = 0.9689654199776964
$ws.Range("O15").Value = 0.7732779360092192
$ws.Range("P15").Value = 0.7732779360092191
$ws.Range("Q15").Value = 40942.30764426383
$ws.Range("R15").Value = 368480.7687983745
$ws.Range("S15").Value = 0.7492795800246593
$ws.Range("T15").Value = 0.7492795800246592
$ws.Range("G16").Value = 4121.876464666667
$ws.Range("H16").Value = 12365.629394
$ws.Range("I16").Value = 0.9689654199776964
$ws.Range("J16").Value = 0.9689654199776964
$ws.Range("M16").Value = 2.334238666666666
$ws.Range("N16").Value = 7.002715999999999
$ws.Range("O16").Value = 0.1817203362411497
$ws.Range("P16").Value = 0.1817203362411496
$ws.Range("Q16").Value = 9621.443423048233
$ws.Range("R16").Value = 86592.99080743409
$ws.Range("S16").Value = 0.1760807219243938
$ws.Range("T16").Value = 0.1760807219243938
$ws.Range("G17").Value = 4121.876464666667
$ws.Range("H17").Value = 12365.629394
$ws.Range("I17").Value = 0.9689654199776964
$ws.Range("J17").Value = 0.9689654199776964
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.07786066666666666
$ws.Range("N17").Value = 0.233582
$ws.Range("O17").Value = 0.006061448098120818
$ws.Range("P17").Value = 0.006061448098120817
$ws.Range("Q17").Value = 320.9320494565898
$ws.Range("R17").Value = 2888.388445109308
$ws.Range("S17").Value = 0.005873333602068647
$ws.Range("T17").Value = 0.005873333602068646
